$d = $word.ActiveDocument

# 1. Replace the placeholder (wrong) login credentials with the corrected,
#    real ones: "are chef/chef." -> "are Administrator/Cod3Can!."
$d.Content.Find.Execute("chef/chef", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Administrator/Cod3Can!", 2)

# 2. Word's automatic "_GoBack" bookmark used to sit between "fil" and
#    "e." in the "...to shorten the URL to the gist file." sentence
#    (the location of the previous edit). After this edit it should sit
#    right after the newly typed credentials, i.e. between "Cod3Can!"
#    and the final "." of the Note sentence. Move it there.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$found = $d.Content
$found.Find.Execute("Administrator/Cod3Can!", $true, $false, $false, $false, `
                     $false, $true, 1, $false, "", 0)

$bmRange = $d.Range($found.End, $found.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
